$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header C1 from "Archivo_Recibido" to "Archivo Recibido"
$ws.Range("C1").Value = "Archivo Recibido"

# Clear the placeholder "NA" values, leaving the cells blank
$ws.Range("C7").Value = ""
$ws.Range("B8:G8").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("C12").Value = ""
$ws.Range("C13").Value = ""
